$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row with another repo name ("production") using the same
# Data Classification / Target Organization values as the other
# "Confidential" repos.
$ws.Range("A9").Value = "production"
$ws.Range("B9").Value = "Confidential"
$ws.Range("C9").Value = "gk-aks-Confidential"

# Update the active selection to reflect the new last-used cell.
$ws.Range("C9").Select()
